# Update the "Price" (D) and "Volume(1h)" (E) columns with freshly scraped
# values, as produced by the "Updated symbol list" GitHub Actions job.
#
# All values in these columns are stored as plain text (they look like
# numbers/percentages, e.g. "305.65" or "0.16%"), so assigning them
# directly via .Value would make Excel auto-convert them into real
# numbers/percentages. To keep them as literal text we briefly switch the
# cell to a text number format while assigning the value, then restore the
# default "Normal" style so the cell's style stays unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$updates = @(
    @{ Row = 2;  D = "305.65";      E = "0.16%" },
    @{ Row = 3;  E = "-0.89%" },
    @{ Row = 4;  D = "5.042";       E = "-0.98%" },
    @{ Row = 5;  D = "0.08049";     E = "-0.57%" },
    @{ Row = 6;  D = "1.905";       E = "-1.57%" },
    @{ Row = 7;  D = "4.155";       E = "-0.73%" },
    @{ Row = 8;  D = "7.781";       E = "0.36%" },
    @{ Row = 9;  D = "0.9214";      E = "-0.68%" },
    @{ Row = 10; D = "0.1281";      E = "-5.98%" },
    @{ Row = 11; D = "0.1916";      E = "-0.06%" },
    @{ Row = 12; D = "0.09065";     E = "-1.47%" },
    @{ Row = 13; D = "0.03449";     E = "1.34%" },
    @{ Row = 14; D = "0.09850";     E = "0.22%" },
    @{ Row = 15; D = "0.001417";    E = "0.71%" },
    @{ Row = 16; D = "0.006228";    E = "7.76%" },
    @{ Row = 17; D = "3.762";       E = "5.30%" },
    @{ Row = 18; D = "3.379";       E = "13.30%" },
    @{ Row = 19; E = "-0.18%" },
    @{ Row = 20; D = "0.1348";      E = "1.23%" },
    @{ Row = 21; D = "5.167";       E = "5.49%" },
    @{ Row = 22; E = "0.23%" },
    @{ Row = 23; D = "0.04430";     E = "0.20%" },
    @{ Row = 24; D = "0.001235";    E = "1.05%" },
    @{ Row = 25; D = "0.004617";    E = "-4.25%" },
    @{ Row = 27; D = "0.0001252";   E = "-3.83%" },
    @{ Row = 28; D = "0.0004448";   E = "42.12%" },
    @{ Row = 39; D = "0.01944";     E = "-3.68%" },
    @{ Row = 40; D = "0.05468";     E = "11.05%" },
    @{ Row = 41; D = "0.007648";    E = "0.54%" },
    @{ Row = 42; D = "0.01012";     E = "-0.87%" },
    @{ Row = 43; E = "-1.76%" },
    @{ Row = 44; D = "0.002153";    E = "2.40%" },
    @{ Row = 45; D = "0.009836";    E = "-10.57%" },
    @{ Row = 46; D = "0.00006130";  E = "-4.20%" },
    @{ Row = 47; E = "0.01%" },
    @{ Row = 48; D = "64.96";       E = "2.19%" },
    @{ Row = 49; D = "0.001660";    E = "39.43%" },
    @{ Row = 50; D = "0.00002103";  E = "0.01%" },
    @{ Row = 51; D = "0.0002003";   E = "0.01%" }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        Set-TextValue ($ws.Cells.Item($u.Row, 4)) $u.D
    }
    if ($u.ContainsKey("E")) {
        Set-TextValue ($ws.Cells.Item($u.Row, 5)) $u.E
    }
}
